$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.760.85"
$ws.Range("D3").Value = "2.252.21"
$ws.Range("E3").Value = "  +3.73%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.36"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.62%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.661"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +17.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.82%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0964"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.21%  "
$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.73%  "
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "2.581.34"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.880"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.27%  "
$ws.Range("D18").Value = "2.251.40"
$ws.Range("E18").Value = "  +1.98%  "
$ws.Range("D19").Value = "42.643.22"
$ws.Range("E19").Value = "  +4.29%  "
$ws.Range("D20").Value = "0.0₃0985"
$ws.Range("E20").Value = "  +4.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.126"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "31.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +21.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0785"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.126"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0319"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.200"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("E48").Value = "  +4.12%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.41%  "
